# Auto-generated edit script: update Leve profit-calculation cells
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed market-board price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 4791.8945
$ws.Cells.Item(43, 9).Value = 3304.7
$ws.Cells.Item(43, 10).Value = 6444.3335
$ws.Cells.Item(43, 11).Value = 3304.7
$ws.Cells.Item(43, 12).Value = 6444.3335
$ws.Cells.Item(43, 13).Value = -3235.7
$ws.Cells.Item(43, 14).Value = -6582.3335
$ws.Cells.Item(138, 8).Value = 2435.5227
$ws.Cells.Item(138, 9).Value = 1913.0385
$ws.Cells.Item(138, 10).Value = 3190.2222
$ws.Cells.Item(138, 11).Value = 5739.1155
$ws.Cells.Item(138, 12).Value = 9570.6666
$ws.Cells.Item(138, 13).Value = -599.1154999999999
$ws.Cells.Item(138, 14).Value = -19850.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(56, 8).Value = 27000
$ws.Cells.Item(56, 9).Value = 20000
$ws.Cells.Item(56, 11).Value = 20000
$ws.Cells.Item(56, 13).Value = -19258

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2000
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 2000
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 2000
$ws.Cells.Item(20, 13).Value = ""
$ws.Cells.Item(20, 14).Value = -2494
$ws.Cells.Item(86, 8).Value = 2043.4286
$ws.Cells.Item(86, 9).Value = 1969.8462
$ws.Cells.Item(86, 10).Value = 3000
$ws.Cells.Item(86, 11).Value = 1969.8462
$ws.Cells.Item(86, 12).Value = 3000
$ws.Cells.Item(86, 13).Value = -846.8462
$ws.Cells.Item(86, 14).Value = -5246
$ws.Cells.Item(89, 8).Value = 2043.4286
$ws.Cells.Item(89, 9).Value = 1969.8462
$ws.Cells.Item(89, 10).Value = 3000
$ws.Cells.Item(89, 11).Value = 9849.231
$ws.Cells.Item(89, 12).Value = 15000
$ws.Cells.Item(89, 13).Value = -4233.231
$ws.Cells.Item(89, 14).Value = -26232
$ws.Cells.Item(107, 8).Value = 1456.6842
$ws.Cells.Item(107, 9).Value = 1398.7646
$ws.Cells.Item(107, 11).Value = 1398.7646
$ws.Cells.Item(107, 13).Value = 521.2354
$ws.Cells.Item(109, 8).Value = 67500
$ws.Cells.Item(109, 9).Value = 60000
$ws.Cells.Item(109, 10).Value = 70000
$ws.Cells.Item(109, 11).Value = 60000
$ws.Cells.Item(109, 12).Value = 70000
$ws.Cells.Item(109, 13).Value = -58613
$ws.Cells.Item(109, 14).Value = -72774
$ws.Cells.Item(134, 8).Value = 1583.381
$ws.Cells.Item(134, 9).Value = 1223.7894
$ws.Cells.Item(134, 11).Value = 3671.3682
$ws.Cells.Item(134, 13).Value = -1136.3682

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 34983
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 34983
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 34983
$ws.Cells.Item(3, 13).Value = ""
$ws.Cells.Item(3, 14).Value = -35209
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 13).Value = ""
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 13).Value = ""
$ws.Cells.Item(99, 8).Value = 2769
$ws.Cells.Item(99, 9).Value = 1913.1428
$ws.Cells.Item(99, 10).Value = 5764.5
$ws.Cells.Item(99, 11).Value = 1913.1428
$ws.Cells.Item(99, 12).Value = 5764.5
$ws.Cells.Item(99, 13).Value = -415.1428000000001
$ws.Cells.Item(99, 14).Value = -8760.5
$ws.Cells.Item(107, 8).Value = 41667070
$ws.Cells.Item(107, 9).Value = 55555816
$ws.Cells.Item(107, 10).Value = 830.3333
$ws.Cells.Item(107, 11).Value = 55555816
$ws.Cells.Item(107, 12).Value = 830.3333
$ws.Cells.Item(107, 13).Value = -55553896
$ws.Cells.Item(107, 14).Value = -4670.3333
$ws.Cells.Item(126, 8).Value = 2769
$ws.Cells.Item(126, 9).Value = 1913.1428
$ws.Cells.Item(126, 10).Value = 5764.5
$ws.Cells.Item(126, 11).Value = 5739.428400000001
$ws.Cells.Item(126, 12).Value = 17293.5
$ws.Cells.Item(126, 13).Value = -3269.428400000001
$ws.Cells.Item(126, 14).Value = -22233.5
$ws.Cells.Item(134, 8).Value = 2189.4092
$ws.Cells.Item(134, 9).Value = 1870.5555
$ws.Cells.Item(134, 11).Value = 5611.666499999999
$ws.Cells.Item(134, 13).Value = -3076.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 2000
$ws.Cells.Item(22, 9).Value = 2000
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 6000
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -5831
$ws.Cells.Item(22, 14).Value = ""
$ws.Cells.Item(25, 8).Value = 333.33334
$ws.Cells.Item(25, 9).Value = 150
$ws.Cells.Item(25, 10).Value = 425
$ws.Cells.Item(25, 11).Value = 450
$ws.Cells.Item(25, 12).Value = 1275
$ws.Cells.Item(25, 13).Value = -281
$ws.Cells.Item(25, 14).Value = -1613
$ws.Cells.Item(27, 8).Value = 2000
$ws.Cells.Item(27, 9).Value = 2000
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 6000
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = -5898
$ws.Cells.Item(27, 14).Value = ""
$ws.Cells.Item(30, 8).Value = 333.33334
$ws.Cells.Item(30, 9).Value = 150
$ws.Cells.Item(30, 10).Value = 425
$ws.Cells.Item(30, 11).Value = 450
$ws.Cells.Item(30, 12).Value = 1275
$ws.Cells.Item(30, 13).Value = -348
$ws.Cells.Item(30, 14).Value = -1479
$ws.Cells.Item(36, 8).Value = 1412.8
$ws.Cells.Item(36, 10).Value = 542.6667
$ws.Cells.Item(36, 12).Value = 1628.0001
$ws.Cells.Item(36, 14).Value = -1966.0001
$ws.Cells.Item(40, 8).Value = 137.6
$ws.Cells.Item(40, 9).Value = 137.6
$ws.Cells.Item(40, 11).Value = 550.4
$ws.Cells.Item(40, 13).Value = -481.4
$ws.Cells.Item(57, 8).Value = 3500
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 13).Value = ""
$ws.Cells.Item(69, 8).Value = 2979.8
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 13).Value = ""
$ws.Cells.Item(72, 8).Value = 2979.8
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 13).Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 766.6667
$ws.Cells.Item(97, 9).Value = 381.81818
$ws.Cells.Item(97, 11).Value = 381.81818
$ws.Cells.Item(97, 13).Value = 114.18182
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 14).Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(5, 8).Value = 15011
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 15011
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 15011
$ws.Cells.Item(5, 13).Value = ""
$ws.Cells.Item(5, 14).Value = -15237
$ws.Cells.Item(23, 8).Value = 7003002
$ws.Cells.Item(23, 9).Value = 7003002
$ws.Cells.Item(23, 11).Value = 7003002
$ws.Cells.Item(23, 13).Value = -7002772
$ws.Cells.Item(34, 8).Value = 19000
$ws.Cells.Item(34, 9).Value = 30000
$ws.Cells.Item(34, 11).Value = 30000
$ws.Cells.Item(34, 13).Value = -29828
$ws.Cells.Item(41, 8).Value = 24466.666
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 24466.666
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 24466.666
$ws.Cells.Item(41, 13).Value = ""
$ws.Cells.Item(41, 14).Value = -25342.666
$ws.Cells.Item(136, 8).Value = 6883
$ws.Cells.Item(136, 9).Value = 6715.9165
$ws.Cells.Item(136, 11).Value = 20147.7495
$ws.Cells.Item(136, 13).Value = -17597.7495

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 19885
$ws.Cells.Item(11, 10).Value = 19885
$ws.Cells.Item(11, 12).Value = 19885
$ws.Cells.Item(11, 14).Value = -20169
$ws.Cells.Item(30, 8).Value = 21505
$ws.Cells.Item(30, 10).Value = 40010
$ws.Cells.Item(30, 12).Value = 40010
$ws.Cells.Item(30, 14).Value = -40224
